$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -3
    5  = -1
    6  = -2
    7  = -4
    8  = 5
    9  = 1
    10 = -2
    11 = 4
    12 = 1
    14 = -2
    15 = -2
    16 = -1
    17 = -2
    18 = 1
    19 = 2
    20 = -5
    22 = 1
    24 = 5
    27 = -1
    30 = 1
    31 = 1
    33 = -2
    34 = 4
    35 = 4
    36 = 2
    37 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
